$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell from "DEMO NUMBER" to "NSS NUMBER"
$ws.Range("A1").Value = "NSS NUMBER"

# Reset the view back to the top-left (A1) and select A1,
# matching the saved view state after the edit.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
